$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, matching the style used by the
# other header cells (e.g. G1: bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the "Save" values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
